$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Update address of the consumer (N2)
$ws.Range("N2").Value = "г. Иркутск, ул. Лермонтова, д. 10, кв. 1"

# Update date of last received payment (O14)
$ws.Range("O14").Value = "15.06.2021г."

# Update consumer's full name (D6)
$ws.Range("D6").Value = "Иванов Иван Иванович"

# Update billed amounts for Электроэнергия (Electricity) row
$ws.Range("E35").Value = 450.0
$ws.Range("G35").Value = 600.0
